$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Harbhajan Singh"

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2 - 10th match vs RCB
$ws.Range("A2").Value = "10th"
$ws.Range("B2").Value = "Kolkata Knight Riders"
$ws.Range("C2").Value = "Harbhajan Singh"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "'2"
$ws.Range("F2").Value = "'2"
$ws.Range("G2").Value = "'0"
$ws.Range("H2").Value = "'0"
$ws.Range("I2").Value = "'100.00"
$ws.Range("J2").Value = "Royal Challengers Bangalore"
$ws.Range("K2").Value = "Chennai"
$ws.Range("L2").Value = "April 18"
$ws.Range("M2").Value = "RCB won by 38 runs"

# Row 3 - 5th match vs Mumbai Indians
$ws.Range("A3").Value = "5th"
$ws.Range("B3").Value = "Kolkata Knight Riders"
$ws.Range("C3").Value = "Harbhajan Singh"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "'2"
$ws.Range("F3").Value = "'2"
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'0"
$ws.Range("I3").Value = "'100.00"
$ws.Range("J3").Value = "Mumbai Indians"
$ws.Range("K3").Value = "Chennai"
$ws.Range("L3").Value = "April 13"
$ws.Range("M3").Value = "Mumbai won by 10 runs"
